$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 751
$ws.Range("I2").Value = 1966
$ws.Range("J2").Value = 7991
$ws.Range("K2").Value = 48
$ws.Range("L2").Value = 2168
$ws.Range("M2").Value = 143
$ws.Range("N2").Value = 1383
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 32
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 111
$ws.Range("S2").Value = 828
$ws.Range("T2").Value = 1383
$ws.Range("U2").Value = 113
$ws.Range("V2").Value = 12655
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 12616
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 175
$ws.Range("AA2").Value = 91
